$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: update "last status check" timestamp from 08:30 to 08:45
$ws.Range("F1").Value = "Last status check on: 08.02.2022 08:45"

# Row 8 (Benzina Albert Modrice): refresh scraped price values
$ws.Range("B8").Value = 36.9
$ws.Range("C8").Value = 36.5

# D8 / E8 switch from numeric (price-delta / serial-date) to plain text,
# matching the new scraper output format. Force text storage (no leftover
# numeric style) by flipping to a text number format, writing the value,
# then resetting the style back to Normal so no stray custom style lingers.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "+0.4"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2022-02-08 08:45:57"
$ws.Range("E8").Style = "Normal"
